$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoan_Input")
$wsOut = $wb.Worksheets.Item("ProductLoan_Output")

# ---------------------------------------------------------------------------
# Rename the product (shared across both sheets' header row, B1)
# ---------------------------------------------------------------------------
$newProductName = "437-RBI-EI-DB-SAR-REC-NON-RNI-CTPD-SAR-MD-TR-1-EarlyRepayment"
$ws.Range("B1").Value = $newProductName
$wsOut.Range("B1").Value = $newProductName

# ---------------------------------------------------------------------------
# shortname (B3): was text "kar9" -> now the plain numeric code 437
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = 437

# ---------------------------------------------------------------------------
# nominalinterestratedefault (B11): 12 -> 1
# ---------------------------------------------------------------------------
$ws.Range("B11").Value = 1

# ---------------------------------------------------------------------------
# Insert new row 22: preclosureinterestcalculationrule
# (Insert() already clones formatting from the row above, matching A21/B21)
# ---------------------------------------------------------------------------
$ws.Rows("22").Insert()
$ws.Range("A22").Value = "preclosureinterestcalculationrule"
$ws.Range("B22").Value = "Calculate till pre closure date"

# ---------------------------------------------------------------------------
# Append new loan-accounting rows 29-40
# ---------------------------------------------------------------------------
$ws.Range("A10:B10").Copy()
$ws.Range("A29:B40").PasteSpecial(-4122)

$accounts = @(
    @("fundsource", "Cash"),
    @("loanprotfolio", "Loan portfolio "),
    @("interestreceivable", "Interest Receivable "),
    @("penaltiesreceivable", "Penalties Receivable "),
    @("transferinsuspense", "Transfer in Suspence "),
    @("feesreceivable", "Fees Receivable"),
    @("incomefrominterest", "Income from interest"),
    @("incomefrompenalties", "Income from penalties"),
    @("incomefromfees", "Income from fees"),
    @("incomefromrecoveryrepayments", "Income from recovery repayments"),
    @("loseswrittenoff", "Losses Writtenoff "),
    @("overpaymentliability", "Overpayment Liability")
)

$r = 29
foreach ($pair in $accounts) {
    $ws.Range("A$r").Value = $pair[0]
    $ws.Range("B$r").Value = $pair[1]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Selection / active cell bookkeeping
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("B1").Select()

$wsOut.Range("B1").Select()
